# Ninez-YA ya4-4-calidad-aire.xlsx — add "metadatos" sheet describing the
# variables in "datos" (renamed from "Sheet1"), and move the active
# selection off of the first sheet onto the new metadata sheet.

$wb = $excel.ActiveWorkbook

# --- rename the original sheet, add the metadata sheet right after it ---
$datos = $wb.Worksheets.Item(1)
$datos.Name = "datos"

$meta = $wb.Worksheets.Add($null, $datos)
$meta.Name = "metadatos"

# --- header row ---
$meta.Range("A1").Value = "Variables"
$meta.Range("B1").Value = "Descripción"
$meta.Range("C1").Value = "Fuente"
$meta.Range("D1").Value = "Fecha_de_extracción"

$fuente = "Sistema Único de Información de Servicios Públicos - SUI"
$fechaExtraccion = 45694   # 2025-02-06, as a raw serial so no ad-hoc number
                            # format gets minted when the value is written

# --- codmpio ---
$meta.Range("A2").Value = "codmpio"
$meta.Range("B2").Value = "Código del municipio"
$meta.Range("C2").Value = $fuente
$meta.Range("D2").NumberFormat = "mm-dd-yy"
$meta.Range("D2").Value = $fechaExtraccion

# --- anno ---
$meta.Range("A3").Value = "anno"
$meta.Range("B3").Value = "Año"
$meta.Range("C3").Value = $fuente
$meta.Range("D3").Value = $fechaExtraccion

# --- concentracion (note: source data has "Fuente" text duplicated in col A) ---
$meta.Range("A4").Font.Name = "Recursive"
$meta.Range("A4").Font.Size = 10
$meta.Range("A4").Font.Color = 5119745
$meta.Range("A4").Value = $fuente
$meta.Range("B4").Value = "Medida Adimensional Calidad de Aire"
$meta.Range("C4").Value = $fuente
$meta.Range("D4").Value = $fechaExtraccion

# re-use the exact same date style (xf index) on D3/D4 instead of minting a
# fresh (but equivalent) style entry for each cell
$meta.Range("D2").Copy()
$meta.Range("D3:D4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- selection: datos -> C1, metadatos (active tab) -> E4 ---
$datos.Range("C1").Select() | Out-Null
$meta.Range("E4").Select() | Out-Null
